$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fundName = "EQUITAS SELECTION"
$assets = @("Stocks", "LFT", "CLCD16", "PETR4", "Compromissadas")

$row = 17
foreach ($asset in $assets) {
    $ws.Cells.Item($row, 1).Value = $fundName
    $ws.Cells.Item($row, 2).Value = $asset
    $row = $row + 1
}

$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()

$ws.Range("A18").Select()
